# edit.ps1
#
# Refresh the cryptocurrency market snapshot (Price, Volume(1h), Hora columns)
# with the values pulled from the latest GitHub Actions data run.
# Column B (Coin), C (Link), F (Data/date) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $text) {
    # The sheet stores these values as text (e.g. "5.030", "-0.81%"), so force
    # a text number format before assigning the value. Otherwise Excel would
    # coerce numeric-looking strings (and drop significant trailing zeros).
    $range.NumberFormat = "@"
    $range.Value = $text
}

# Row 2
Set-TextCell $ws.Range("D2") "306.27"
Set-TextCell $ws.Range("E2") "-0.38%"
Set-TextCell $ws.Range("G2") "19"

# Row 3
Set-TextCell $ws.Range("D3") "36.29"
Set-TextCell $ws.Range("E3") "-0.93%"
Set-TextCell $ws.Range("G3") "19"

# Row 4
Set-TextCell $ws.Range("D4") "5.039"
Set-TextCell $ws.Range("E4") "0.09%"
Set-TextCell $ws.Range("G4") "19"

# Row 5
Set-TextCell $ws.Range("D5") "0.07892"
Set-TextCell $ws.Range("E5") "0.39%"
Set-TextCell $ws.Range("G5") "19"

# Row 6
Set-TextCell $ws.Range("D6") "2.138"
Set-TextCell $ws.Range("E6") "-2.04%"
Set-TextCell $ws.Range("G6") "19"

# Row 7
Set-TextCell $ws.Range("D7") "7.972"
Set-TextCell $ws.Range("E7") "-1.09%"
Set-TextCell $ws.Range("G7") "19"

# Row 8
Set-TextCell $ws.Range("D8") "0.9283"
Set-TextCell $ws.Range("E8") "0.38%"
Set-TextCell $ws.Range("G8") "19"

# Row 9
Set-TextCell $ws.Range("D9") "0.09689"
Set-TextCell $ws.Range("E9") "-2.56%"
Set-TextCell $ws.Range("G9") "19"

# Row 10
Set-TextCell $ws.Range("D10") "0.1865"
Set-TextCell $ws.Range("E10") "-0.80%"
Set-TextCell $ws.Range("G10") "19"

# Row 11
Set-TextCell $ws.Range("D11") "0.09023"
Set-TextCell $ws.Range("E11") "3.74%"
Set-TextCell $ws.Range("G11") "19"

# Row 12
Set-TextCell $ws.Range("D12") "0.03693"
Set-TextCell $ws.Range("E12") "2.03%"
Set-TextCell $ws.Range("G12") "19"

# Row 13
Set-TextCell $ws.Range("D13") "0.09888"
Set-TextCell $ws.Range("E13") "-0.60%"
Set-TextCell $ws.Range("G13") "19"

# Row 14
Set-TextCell $ws.Range("D14") "0.001433"
Set-TextCell $ws.Range("E14") "-3.60%"
Set-TextCell $ws.Range("G14") "19"

# Row 15
Set-TextCell $ws.Range("D15") "0.005628"
Set-TextCell $ws.Range("E15") "-0.45%"
Set-TextCell $ws.Range("G15") "19"

# Row 16
Set-TextCell $ws.Range("D16") "3.484"
Set-TextCell $ws.Range("E16") "0.67%"
Set-TextCell $ws.Range("G16") "19"

# Row 17
Set-TextCell $ws.Range("D17") "4.166"
Set-TextCell $ws.Range("E17") "2.59%"
Set-TextCell $ws.Range("G17") "19"

# Row 18
Set-TextCell $ws.Range("E18") "19.72%"
Set-TextCell $ws.Range("G18") "19"

# Row 19
Set-TextCell $ws.Range("E19") "-0.83%"
Set-TextCell $ws.Range("G19") "19"

# Row 20
Set-TextCell $ws.Range("D20") "0.1334"
Set-TextCell $ws.Range("E20") "-0.25%"
Set-TextCell $ws.Range("G20") "19"

# Row 21
Set-TextCell $ws.Range("D21") "5.081"
Set-TextCell $ws.Range("E21") "3.03%"
Set-TextCell $ws.Range("G21") "19"

# Row 22
Set-TextCell $ws.Range("D22") "0.2254"
Set-TextCell $ws.Range("E22") "2.32%"
Set-TextCell $ws.Range("G22") "19"

# Row 23
Set-TextCell $ws.Range("D23") "0.04583"
Set-TextCell $ws.Range("E23") "-0.84%"
Set-TextCell $ws.Range("G23") "19"

# Row 24
Set-TextCell $ws.Range("D24") "0.001237"
Set-TextCell $ws.Range("E24") "0.25%"
Set-TextCell $ws.Range("G24") "19"

# Row 25
Set-TextCell $ws.Range("D25") "0.004797"
Set-TextCell $ws.Range("E25") "-7.76%"
Set-TextCell $ws.Range("G25") "19"

# Row 26
Set-TextCell $ws.Range("D26") "0.0001303"
Set-TextCell $ws.Range("E26") "-7.04%"
Set-TextCell $ws.Range("G26") "19"

# Row 27
Set-TextCell $ws.Range("E27") "74.14%"
Set-TextCell $ws.Range("G27") "19"

# Row 28
Set-TextCell $ws.Range("G28") "19"

# Row 29
Set-TextCell $ws.Range("G29") "19"

# Row 30
Set-TextCell $ws.Range("G30") "19"

# Row 31
Set-TextCell $ws.Range("G31") "19"

# Row 32
Set-TextCell $ws.Range("G32") "19"

# Row 33
Set-TextCell $ws.Range("G33") "19"

# Row 34
Set-TextCell $ws.Range("G34") "19"

# Row 35
Set-TextCell $ws.Range("G35") "19"

# Row 36
Set-TextCell $ws.Range("G36") "19"

# Row 37
Set-TextCell $ws.Range("G37") "19"

# Row 38
Set-TextCell $ws.Range("G38") "19"

# Row 39
Set-TextCell $ws.Range("D39") "0.01985"
Set-TextCell $ws.Range("E39") "9.15%"
Set-TextCell $ws.Range("G39") "19"

# Row 40
Set-TextCell $ws.Range("D40") "0.04907"
Set-TextCell $ws.Range("E40") "3.17%"
Set-TextCell $ws.Range("G40") "19"

# Row 41
Set-TextCell $ws.Range("D41") "0.007840"
Set-TextCell $ws.Range("E41") "-0.53%"
Set-TextCell $ws.Range("G41") "19"

# Row 42
Set-TextCell $ws.Range("D42") "0.1392"
Set-TextCell $ws.Range("E42") "-1.36%"
Set-TextCell $ws.Range("G42") "19"

# Row 43
Set-TextCell $ws.Range("D43") "0.007820"
Set-TextCell $ws.Range("E43") "2.80%"
Set-TextCell $ws.Range("G43") "19"

# Row 44
Set-TextCell $ws.Range("D44") "0.002147"
Set-TextCell $ws.Range("E44") "-1.62%"
Set-TextCell $ws.Range("G44") "19"

# Row 45
Set-TextCell $ws.Range("D45") "0.01126"
Set-TextCell $ws.Range("E45") "11.38%"
Set-TextCell $ws.Range("G45") "19"

# Row 46
Set-TextCell $ws.Range("D46") "0.00006270"
Set-TextCell $ws.Range("E46") "-0.75%"
Set-TextCell $ws.Range("G46") "19"

# Row 47
Set-TextCell $ws.Range("D47") "0.00000000752"
Set-TextCell $ws.Range("E47") "0.10%"
Set-TextCell $ws.Range("G47") "19"

# Row 48
Set-TextCell $ws.Range("E48") "0.06%"
Set-TextCell $ws.Range("G48") "19"

# Row 49
Set-TextCell $ws.Range("D49") "51.71"
Set-TextCell $ws.Range("E49") "57.83%"
Set-TextCell $ws.Range("G49") "19"

# Row 50
Set-TextCell $ws.Range("D50") "0.001904"
Set-TextCell $ws.Range("E50") "-29.30%"
Set-TextCell $ws.Range("G50") "19"

# Row 51
Set-TextCell $ws.Range("D51") "0.00002105"
Set-TextCell $ws.Range("E51") "0.10%"
Set-TextCell $ws.Range("G51") "19"
